# Add season-record columns (Wins / Losses / Ties) to the stats table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting of the existing header cell (bold font, border,
# centered alignment) onto the new header cells so they match the rest
# of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows --------------------------------------------------------------
# Every player on the roster shares the same team season record.
$ws.Range("AD2:AD66").Value = 67
$ws.Range("AE2:AE66").Value = 95
$ws.Range("AF2:AF66").Value = 0
